$d = $word.ActiveDocument

# Locate the paragraph that starts the "Weekly report 4/14/19" entry
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Contains("Weekly report 4/14/19")) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    Write-Host "ERROR: could not locate target paragraph"
} else {
    $targetPara = $d.Paragraphs.Item($targetIndex)
    $targetRange = $targetPara.Range

    $xmlPayload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00D60735" w:rsidRDefault="00D60735" w:rsidP="00D60735"><w:r><w:t xml:space="preserve">Weekly report </w:t></w:r><w:r><w:t xml:space="preserve">for week ending </w:t></w:r><w:r><w:t xml:space="preserve">4/14/19 Performed code review on </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Tuedsay</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> as the author of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>saveData</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) class, and </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>IDataIO</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> interface with applicable main function calls. Identified problems with the code. After the review, changes were made to improve the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>saveData</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>class(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) and the class now writes the test.xml file correctly. Performed code review as moderator on Thursday with functions in the main. </w:t></w:r><w:r><w:t>Thursday was in class code review and on Sunday assisted in debugging the load xml function.</w:t></w:r></w:p><w:p><w:r><w:lastRenderedPageBreak/><w:t>Weekly report for week ending 4</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t xml:space="preserve">/21/19 – Finished the iteration which required loading and saving data from a file. Refactored the </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r><w:t>savexml</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">) function to assist in completion of the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>loadXML</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>() function.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $targetRange.InsertXML($xmlPayload)

    # The original paragraph's content has been pushed down into a now-empty
    # leftover paragraph (InsertXML inserts the new paragraphs *before* the
    # range and empties out the original). Remove that leftover paragraph,
    # including its paragraph mark, so it doesn't leave a stray blank line.
    $leftover = $d.Paragraphs.Item($targetIndex + 2)
    $cleanupRange = $d.Range($leftover.Range.Start - 1, $leftover.Range.End)
    $cleanupRange.Delete()
}
